$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1822222222222222
$ws.Range("C2").Value = 0.5777777777777777
$ws.Range("J2").Value = 0.01777777777777778
$ws.Range("P2").Value = 0.1377777777777778
$ws.Range("S2").Value = 0.08444444444444445
$ws.Range("B3").Value = 0.01515151515151515
$ws.Range("C3").Value = 0.01515151515151515
$ws.Range("J3").Value = 0.007575757575757576
$ws.Range("P3").Value = 0.7272727272727273
$ws.Range("S3").Value = 0.2348484848484849
$ws.Range("J4").Value = 0.025
$ws.Range("P4").Value = 0.675
$ws.Range("S4").Value = 0.3
$ws.Range("B6").Value = 0.06796116504854369
$ws.Range("D6").Value = 0.01941747572815534
$ws.Range("F6").Value = 0.06796116504854369
$ws.Range("J6").Value = 0.2524271844660194
$ws.Range("O6").Value = 0.01941747572815534
$ws.Range("Q6").Value = 0.1553398058252427
$ws.Range("R6").Value = 0.1019417475728155
$ws.Range("S6").Value = 0.3155339805825243
$ws.Range("B7").Value = 0.1092896174863388
$ws.Range("D7").Value = 0.01639344262295082
$ws.Range("E7").Value = 0.00546448087431694
$ws.Range("F7").Value = 0.04918032786885246
$ws.Range("J7").Value = 0.1530054644808743
$ws.Range("O7").Value = 0.00546448087431694
$ws.Range("Q7").Value = 0.185792349726776
$ws.Range("R7").Value = 0.03825136612021858
$ws.Range("S7").Value = 0.4371584699453552
$ws.Range("B8").Value = 0.07196029776674938
$ws.Range("D8").Value = 0.02481389578163772
$ws.Range("F8").Value = 0.06947890818858561
$ws.Range("J8").Value = 0.1066997518610422
$ws.Range("O8").Value = 0.02233250620347394
$ws.Range("Q8").Value = 0.1662531017369727
$ws.Range("R8").Value = 0.09677419354838709
$ws.Range("S8").Value = 0.4416873449131514
$ws.Range("B9").Value = 0.08917197452229299
$ws.Range("D9").Value = 0.01910828025477707
$ws.Range("F9").Value = 0.03821656050955414
$ws.Range("J9").Value = 0.07643312101910828
$ws.Range("O9").Value = 0.01910828025477707
$ws.Range("Q9").Value = 0.1656050955414013
$ws.Range("R9").Value = 0.09554140127388536
$ws.Range("S9").Value = 0.4968152866242038
$ws.Range("B10").Value = 0.09058295964125561
$ws.Range("D10").Value = 0.02152466367713005
$ws.Range("E10").Value = 0.0008968609865470852
$ws.Range("F10").Value = 0.08071748878923767
$ws.Range("J10").Value = 0.1022421524663677
$ws.Range("O10").Value = 0.01524663677130045
$ws.Range("Q10").Value = 0.2242152466367713
$ws.Range("R10").Value = 0.08430493273542601
$ws.Range("S10").Value = 0.3802690582959641
$ws.Range("G11").Value = 0.1230283911671924
$ws.Range("J11").Value = 0.110410094637224
$ws.Range("K11").Value = 0.2050473186119874
$ws.Range("L11").Value = 0.555205047318612
$ws.Range("S11").Value = 0.006309148264984227
$ws.Range("G12").Value = 0.7065217391304348
$ws.Range("J12").Value = 0.1793478260869565
$ws.Range("K12").Value = 0.005434782608695652
$ws.Range("L12").Value = 0.03260869565217391
$ws.Range("S12").Value = 0.07608695652173914
$ws.Range("G13").Value = 0.6129032258064516
$ws.Range("J13").Value = 0.3870967741935484
$ws.Range("F15").Value = 0.01731601731601732
$ws.Range("H15").Value = 0.1818181818181818
$ws.Range("I15").Value = 0.09090909090909091
$ws.Range("J15").Value = 0.3246753246753247
$ws.Range("K15").Value = 0.09090909090909091
$ws.Range("M15").Value = 0.008658008658008658
$ws.Range("O15").Value = 0.07792207792207792
$ws.Range("S15").Value = 0.2077922077922078
$ws.Range("F16").Value = 0.02013422818791946
$ws.Range("H16").Value = 0.174496644295302
$ws.Range("I16").Value = 0.06711409395973154
$ws.Range("J16").Value = 0.3758389261744967
$ws.Range("K16").Value = 0.174496644295302
$ws.Range("M16").Value = 0.006711409395973154
$ws.Range("N16").Value = 0.006711409395973154
$ws.Range("O16").Value = 0.08053691275167785
$ws.Range("S16").Value = 0.09395973154362416
$ws.Range("F17").Value = 0.0293398533007335
$ws.Range("H17").Value = 0.2029339853300733
$ws.Range("I17").Value = 0.09290953545232274
$ws.Range("J17").Value = 0.3716381418092909
$ws.Range("K17").Value = 0.09046454767726161
$ws.Range("M17").Value = 0.01466992665036675
$ws.Range("O17").Value = 0.06845965770171149
$ws.Range("S17").Value = 0.1295843520782396
$ws.Range("F18").Value = 0.01129943502824859
$ws.Range("H18").Value = 0.1807909604519774
$ws.Range("I18").Value = 0.07909604519774012
$ws.Range("J18").Value = 0.4463276836158192
$ws.Range("K18").Value = 0.096045197740113
$ws.Range("M18").Value = 0.01129943502824859
$ws.Range("O18").Value = 0.05649717514124294
$ws.Range("S18").Value = 0.1186440677966102
$ws.Range("F19").Value = 0.01885175664095973
$ws.Range("H19").Value = 0.1910882604970009
$ws.Range("I19").Value = 0.06341045415595545
$ws.Range("J19").Value = 0.3693230505569837
$ws.Range("K19").Value = 0.1259640102827763
$ws.Range("M19").Value = 0.02056555269922879
$ws.Range("N19").Value = 0.002570694087403599
$ws.Range("O19").Value = 0.0831191088260497
$ws.Range("S19").Value = 0.1251071122536418
